$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 61; existing rows 61:85 shift down to 62:86.
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with the new weekly price record.
$ws.Cells.Item(61, 1).Value = 5
$ws.Cells.Item(61, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(61, 3).Value = "Maule"
$ws.Cells.Item(61, 4).Value = 44813
$ws.Cells.Item(61, 5).Value = 7
$ws.Cells.Item(61, 6).Value = 100112026
$ws.Cells.Item(61, 7).Value = "Haba"
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 200
$ws.Cells.Item(61, 11).Value = 12000
$ws.Cells.Item(61, 12).Value = 12000
$ws.Cells.Item(61, 13).Value = 12000
$ws.Cells.Item(61, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(61, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(61, 16).Value = 480
$ws.Cells.Item(61, 17).Value = 25
$ws.Cells.Item(61, 18).Value = "Hortaliza"
